$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "Family services" layout change: add a new Year / Number pair of columns
# (V:W) showing yearly counts, next to the existing Category/value pairs.
$ws.Range("V1").Value = "Year "
$ws.Range("W1").Value = "Number"

$ws.Range("V2").Value = 2016
$ws.Range("W2").Value = 4

$ws.Range("V3").Value = 2017
$ws.Range("W3").Value = 10

$ws.Range("V4").Value = 2018
$ws.Range("W4").Value = 7

$ws.Range("V5").Value = 2019
$ws.Range("W5").Value = 20

$ws.Range("V6").Value = 2020
$ws.Range("W6").Value = 22

# Scroll over to show the newly added columns and select the cell below the
# last entry, matching the reviewed state of the sheet.
$ws.Range("W7").Select()
